$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.3
$ws.Range("I2").Value = 3.1
$ws.Range("L2").Value = 3.5
$ws.Range("U2").Value = 1.57
$ws.Range("V2").Value = 2.25
$ws.Range("AJ2").Value = 11
$ws.Range("AN2").Value = 4.5
$ws.Range("AX2").Value = 15
